$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A13 with consolidated tuple-like strings
$ws.Range("A2").Value = "('Assassin', ['Token Creature — Assassin', 'Whenever this creature deals combat damage to a player, that player loses the game.', '1/1'])"
$ws.Range("A3").Value = "('Bird', ['Token Creature — Bird', 'Flying', '1/1'])"
$ws.Range("A4").Value = "('Centaur', ['Token Creature — Centaur', '3/3'])"
$ws.Range("A5").Value = "('Dragon', ['Token Creature — Dragon', 'Flying', '6/6'])"
$ws.Range("A6").Value = "('Elemental', ['Token Creature — Elemental', 'Vigilance', '8/8'])"
$ws.Range("A7").Value = "('Goblin', ['Token Creature — Goblin', '1/1'])"
$ws.Range("A8").Value = "('Knight', ['Token Creature — Knight', 'Vigilance', '2/2'])"
$ws.Range("A9").Value = "('Ooze', ['Token Creature — Ooze', '*/*'])"
$ws.Range("A10").Value = "('Rhino', ['Token Creature — Rhino', 'Trample', '4/4'])"
$ws.Range("A11").Value = "('Saproling', ['Token Creature — Saproling', '1/1'])"
$ws.Range("A12").Value = "('Soldier', ['Token Creature — Soldier', '1/1'])"
$ws.Range("A13").Value = "('Wurm', ['Token Creature — Wurm', 'Trample', '5/5'])"

# Remove the now-obsolete rows 14 through 44
$ws.Range("A14:A44").EntireRow.Delete()
